# "solved problem about become VIP"
#
# 1) Drink sheet: remove the stray duplicate "Latte" row (row 6) that was
#    left over from debugging the VIP logic.
# 2) Room sheet: correct the first booking record (row 2) — price, booking
#    date and the booking user id were wrong, which is what broke the
#    "become VIP" check. Values are written as plain text (matching how the
#    row is now exported/consumed downstream).

$wb = $excel.ActiveWorkbook

# --- Drink: delete row 6 entirely ---------------------------------------
$drink = $wb.Worksheets.Item("Drink")
$drink.Rows.Item(6).Delete()

# --- Room: fix booking row 2 --------------------------------------------
$room = $wb.Worksheets.Item("Room")

$room.Range("A2").Value = "'1"
$room.Range("C2").Value = "'188.0"
$room.Range("D2").Value = "'2022-12-11"
$room.Range("E2").Value = "'5"
$room.Range("F2").Value = "'10"
